$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# --- Row 4 (1 thread run) ---
# C4 used to be the formula =C6*2 (26897092); it becomes a plain measured value.
$ws.Range("C4").Value = 26690977

# K4/L4 used to derive from K6/L6 (=2*K6 / =2*L6); they now derive straight from C4/K4,
# matching the other rows' pattern (run time in minutes / hours).
$ws.Range("K4").Formula = "=C4/1000/60"
$ws.Range("L4").Formula = "=K4/60"

# New M4 cell: speed-up of the 1-thread run relative to itself (always 1).
$ws.Range("M4").Formula = "=C4/C4"

# --- New row 22 data (SGSIM original baseline run) ---
$ws.Range("C22").Value = 10043951
$ws.Range("K22").Formula = "=C22/1000/60"
$ws.Range("K22").NumberFormat = "0.00"
$ws.Range("L22").Formula = "=K22/60"
$ws.Range("L22").NumberFormat = "0.00"

# --- View state: scroll/selection moved ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B24").Select()
